# Auto-generated edit script: updates computed profit/price columns (H:N)
# across multiple worksheets to match the target snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4767.7896
$ws.Range("I137").Value = 2371.818
$ws.Range("J137").Value = 8062.25
$ws.Range("K137").Value = 7115.454000000001
$ws.Range("L137").Value = 24186.75
$ws.Range("M137").Value = -4565.454000000001
$ws.Range("N137").Value = -29286.75
$ws.Range("H138").Value = 6692.4116
$ws.Range("J138").Value = 5872.4194
$ws.Range("L138").Value = 17617.2582
$ws.Range("N138").Value = -27897.2582

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2635.8276
$ws.Range("J2").Value = 7775.5713
$ws.Range("L2").Value = 7775.5713
$ws.Range("N2").Value = -8001.5713
$ws.Range("H97").Value = 3552.25
$ws.Range("I97").Value = 3236.3333
$ws.Range("K97").Value = 3236.3333
$ws.Range("M97").Value = -2740.3333
$ws.Range("H116").Value = 2635.8276
$ws.Range("J116").Value = 7775.5713
$ws.Range("L116").Value = 7775.5713
$ws.Range("N116").Value = -12363.5713
$ws.Range("H122").Value = 1622.24
$ws.Range("I122").Value = 1600.326
$ws.Range("J122").Value = 1874.25
$ws.Range("K122").Value = 4800.978
$ws.Range("L122").Value = 5622.75
$ws.Range("M122").Value = -2350.978
$ws.Range("N122").Value = -10522.75
$ws.Range("H132").Value = 1728.2727
$ws.Range("I132").Value = 1188.4333
$ws.Range("J132").Value = 2885.0715
$ws.Range("K132").Value = 3565.2999
$ws.Range("L132").Value = 8655.2145
$ws.Range("M132").Value = -1035.2999
$ws.Range("N132").Value = -13715.2145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2635.8276
$ws.Range("J3").Value = 7775.5713
$ws.Range("L3").Value = 7775.5713
$ws.Range("N3").Value = -8003.5713
$ws.Range("H20").Value = 15628242
$ws.Range("I20").Value = 18522098
$ws.Range("J20").Value = 1425.2
$ws.Range("K20").Value = 18522098
$ws.Range("L20").Value = 1425.2
$ws.Range("M20").Value = -18521851
$ws.Range("N20").Value = -1919.2
$ws.Range("H80").Value = 584
$ws.Range("J80").Value = 638.8889
$ws.Range("L80").Value = 638.8889
$ws.Range("N80").Value = -2634.8889
$ws.Range("H83").Value = 584
$ws.Range("J83").Value = 638.8889
$ws.Range("L83").Value = 3194.4445
$ws.Range("N83").Value = -13178.4445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3533.2156
$ws.Range("I31").Value = 1665.6666
$ws.Range("J31").Value = 3649.9375
$ws.Range("K31").Value = 1665.6666
$ws.Range("L31").Value = 3649.9375
$ws.Range("M31").Value = -1370.6666
$ws.Range("N31").Value = -4239.9375
$ws.Range("H34").Value = 3533.2156
$ws.Range("I34").Value = 1665.6666
$ws.Range("J34").Value = 3649.9375
$ws.Range("K34").Value = 1665.6666
$ws.Range("L34").Value = 3649.9375
$ws.Range("M34").Value = -1463.6666
$ws.Range("N34").Value = -4053.9375
$ws.Range("H138").Value = 69914
$ws.Range("J138").Value = 69914
$ws.Range("L138").Value = 69914
$ws.Range("N138").Value = -80194

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H28").Value = 194
$ws.Range("I28").Value = 179
$ws.Range("K28").Value = 537
$ws.Range("M28").Value = -305
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()
$ws.Range("H68").Value = 6353499
$ws.Range("I68").Value = 5557372
$ws.Range("J68").Value = 6671950
$ws.Range("K68").Value = 16672116
$ws.Range("L68").Value = 20015850
$ws.Range("M68").Value = -16671305
$ws.Range("N68").Value = -20017472
$ws.Range("H71").Value = 6353499
$ws.Range("I71").Value = 5557372
$ws.Range("J71").Value = 6671950
$ws.Range("K71").Value = 50016348
$ws.Range("L71").Value = 60047550
$ws.Range("M71").Value = -50012292
$ws.Range("N71").Value = -60055662
$ws.Range("H86").Value = 1821.5
$ws.Range("J86").Value = 2412.8572
$ws.Range("L86").Value = 7238.571599999999
$ws.Range("N86").Value = -9610.571599999999
$ws.Range("H89").Value = 1821.5
$ws.Range("J89").Value = 2412.8572
$ws.Range("L89").Value = 21715.7148
$ws.Range("N89").Value = -33571.7148
$ws.Range("H111").Value = 3500
$ws.Range("I111").Value = 3500
$ws.Range("K111").Value = 10500
$ws.Range("M111").Value = -7433
$ws.Range("H129").Value = 1589.6
$ws.Range("I129").Value = 987
$ws.Range("J129").Value = 4000
$ws.Range("K129").Value = 2961
$ws.Range("L129").Value = 12000
$ws.Range("M129").Value = 2039
$ws.Range("N129").Value = -22000
$ws.Range("H131").Value = 7625.3335
$ws.Range("J131").Value = 2090.4348
$ws.Range("L131").Value = 6271.3044
$ws.Range("N131").Value = -16351.3044

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7399.8
$ws.Range("I70").Value = 4999.5
$ws.Range("K70").Value = 4999.5
$ws.Range("M70").Value = -4729.5
$ws.Range("H73").Value = 7399.8
$ws.Range("I73").Value = 4999.5
$ws.Range("K73").Value = 4999.5
$ws.Range("M73").Value = -4063.5
$ws.Range("H102").Value = 19610488
$ws.Range("I102").Value = 37040524
$ws.Range("K102").Value = 37040524
$ws.Range("M102").Value = -37038902
$ws.Range("H122").Value = 20836192
$ws.Range("I122").Value = 2473.7097
$ws.Range("K122").Value = 7421.1291
$ws.Range("M122").Value = -4971.1291
$ws.Range("H132").Value = 1927.721
$ws.Range("I132").Value = 1463.125
$ws.Range("K132").Value = 4389.375
$ws.Range("M132").Value = -1859.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1970.3158
$ws.Range("I93").Value = 2082.3333
$ws.Range("K93").Value = 2082.3333
$ws.Range("M93").Value = -834.3332999999998
$ws.Range("H122").Value = 4333.8237
$ws.Range("I122").Value = 4161.1816
$ws.Range("K122").Value = 12483.5448
$ws.Range("M122").Value = -10033.5448
$ws.Range("H132").Value = 3144.0232
$ws.Range("I132").Value = 3287.3
$ws.Range("J132").Value = 2813.3845
$ws.Range("K132").Value = 9861.900000000001
$ws.Range("L132").Value = 8440.1535
$ws.Range("M132").Value = -7331.900000000001
$ws.Range("N132").Value = -13500.1535
$ws.Range("H136").Value = 6297.4165
$ws.Range("J136").Value = 6262.3335
$ws.Range("L136").Value = 18787.0005
$ws.Range("N136").Value = -23887.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 508.14285
$ws.Range("I107").Value = 587.6
$ws.Range("J107").Value = 309.5
$ws.Range("K107").Value = 1762.8
$ws.Range("L107").Value = 928.5
$ws.Range("M107").Value = 157.1999999999998
$ws.Range("N107").Value = -4768.5
$ws.Range("H122").Value = 8335762.5
$ws.Range("I122").Value = 2599.3635
$ws.Range("J122").Value = 31251962
$ws.Range("K122").Value = 7798.0905
$ws.Range("L122").Value = 93755886
$ws.Range("M122").Value = -5348.0905
$ws.Range("N122").Value = -93760786
$ws.Range("H132").Value = 1668.2094
$ws.Range("I132").Value = 1428.8572
$ws.Range("J132").Value = 2715.375
$ws.Range("K132").Value = 4286.571599999999
$ws.Range("L132").Value = 8146.125
$ws.Range("M132").Value = -1756.571599999999
$ws.Range("N132").Value = -13206.125
